# Updates cryptos list (Price/Volume(1h) columns, and a row-48/49 coin swap)
# as produced by the scheduled GitHub Actions scraper run.
#
# Price (column D) and Volume(1h) (column E) are stored as plain text in the
# workbook (t="inlineStr"), even when the text looks like a number (e.g.
# "69.416.04", "0.999", "1.00"). Assigning a bare numeric-looking string to
# .Value would make Excel auto-coerce the cell to a real number (losing
# trailing zeros / multi-dot "thousands" formatting such as "69.416.04").
# Prefixing the string with a leading apostrophe forces Excel to treat the
# input as literal text (like typing '69.416.04 into a cell), and the
# subsequent Style reset keeps the cell's style index unchanged (no
# quotePrefix / text-format style lingering on the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'69.416.04"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -1.90%  "
$ws.Cells.Item(3, 4).Value = "'3.686.02"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -2.92%  "
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Value = "'681.50"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -2.92%  "
$ws.Cells.Item(6, 4).Value = "'162.59"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -4.28%  "
$ws.Cells.Item(7, 4).Value = "'3.685.54"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.93%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 5).Value = "  -4.18%  "
$ws.Cells.Item(10, 5).Value = "  -7.55%  "
$ws.Cells.Item(11, 4).Value = "'7.36"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.69%  "
$ws.Cells.Item(12, 4).Value = "'0.446"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -2.52%  "
$ws.Cells.Item(13, 5).Value = "  -4.45%  "
$ws.Cells.Item(14, 4).Value = "'33.56"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -5.88%  "
$ws.Cells.Item(15, 4).Value = "'4.299.09"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -3.16%  "
$ws.Cells.Item(16, 4).Value = "'3.691.13"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -2.54%  "
$ws.Cells.Item(17, 4).Value = "'69.416.25"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -1.91%  "
$ws.Cells.Item(18, 5).Value = "  -0.88%  "
$ws.Cells.Item(19, 5).Value = "  -6.28%  "
$ws.Cells.Item(20, 4).Value = "'6.62"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -6.62%  "
$ws.Cells.Item(21, 4).Value = "'483.85"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.76%  "
$ws.Cells.Item(22, 4).Value = "'9.89"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -6.96%  "
$ws.Cells.Item(23, 5).Value = "  -7.37%  "
$ws.Cells.Item(24, 4).Value = "'80.38"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -4.44%  "
$ws.Cells.Item(25, 4).Value = "'3.831.90"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.97%  "
$ws.Cells.Item(26, 5).Value = "  -8.27%  "
$ws.Cells.Item(27, 5).Value = "  -0.03%  "
$ws.Cells.Item(28, 4).Value = "'11.47"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -4.52%  "
$ws.Cells.Item(29, 5).Value = "  -6.51%  "
$ws.Cells.Item(30, 5).Value = "  -8.29%  "
$ws.Cells.Item(31, 4).Value = "'2.73"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -9.38%  "
$ws.Cells.Item(32, 5).Value = "  -7.25%  "
$ws.Cells.Item(33, 5).Value = "  -6.19%  "
$ws.Cells.Item(34, 4).Value = "'27.07"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -6.43%  "
$ws.Cells.Item(35, 5).Value = "  -4.81%  "
$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.13%  "
$ws.Cells.Item(37, 4).Value = "'3.656.70"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.81%  "
$ws.Cells.Item(38, 4).Value = "'8.51"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -5.68%  "
$ws.Cells.Item(39, 4).Value = "'6.34"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +6.99%  "
$ws.Cells.Item(40, 4).Value = "'0.0936"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -7.16%  "
$ws.Cells.Item(41, 4).Value = "'2.24"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -5.05%  "
$ws.Cells.Item(43, 5).Value = "  -0.04%  "
$ws.Cells.Item(44, 5).Value = "  -6.08%  "
$ws.Cells.Item(45, 4).Value = "'160.16"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -3.88%  "
$ws.Cells.Item(46, 4).Value = "'48.45"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.81%  "
$ws.Cells.Item(47, 4).Value = "'2.86"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -10.93%  "
$ws.Cells.Item(48, 2).Value = "FLOKI"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(48, 4).Value = "'0.000291"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -7.36%  "
$ws.Cells.Item(49, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(49, 4).Value = "'29.87"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +6.51%  "
$ws.Cells.Item(50, 5).Value = "  +1.89%  "
$ws.Cells.Item(51, 4).Value = "'393.48"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -5.75%  "
